$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.554827
$ws.Range("H2").Value = 10.664481
$ws.Range("I2").Value = 0.2148969460055877
$ws.Range("J2").Value = 0.2148969460055877
$ws.Range("O2").Value = 0.02266023449704293
$ws.Range("P2").Value = 0.02266023449704293
$ws.Range("Q2").Value = 0.5597276551696666
$ws.Range("R2").Value = 5.037548896526999
$ws.Range("S2").Value = 0.00486961518918499
$ws.Range("T2").Value = 0.00486961518918499

# Row 3
$ws.Range("G3").Value = 3.554827
$ws.Range("H3").Value = 10.664481
$ws.Range("I3").Value = 0.2148969460055877
$ws.Range("J3").Value = 0.2148969460055877
$ws.Range("M3").Value = 6.739756333333333
$ws.Range("O3").Value = 0.9699521281096917
$ws.Range("P3").Value = 0.9699521281096917
$ws.Range("Q3").Value = 23.95866778715433
$ws.Range("R3").Value = 215.628010084389
$ws.Range("S3").Value = 0.2084397501023933
$ws.Range("T3").Value = 0.2084397501023933

# Row 4
$ws.Range("G4").Value = 3.554827
$ws.Range("H4").Value = 10.664481
$ws.Range("I4").Value = 0.2148969460055877
$ws.Range("J4").Value = 0.2148969460055877
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05133333333333333
$ws.Range("N4").Value = 0.154
$ws.Range("O4").Value = 0.00738763739326543
$ws.Range("P4").Value = 0.00738763739326543
$ws.Range("Q4").Value = 0.1824811193333333
$ws.Range("R4").Value = 1.642330074
$ws.Range("S4").Value = 0.001587580714009422
$ws.Range("T4").Value = 0.001587580714009422

# Row 5
$ws.Range("I5").Value = 0.3107709374420163
$ws.Range("J5").Value = 0.3107709374420163
$ws.Range("O5").Value = 0.02266023449704293
$ws.Range("P5").Value = 0.02266023449704293
$ws.Range("S5").Value = 0.007042142317301948
$ws.Range("T5").Value = 0.007042142317301947

# Row 6
$ws.Range("I6").Value = 0.3107709374420163
$ws.Range("J6").Value = 0.3107709374420163
$ws.Range("M6").Value = 6.739756333333333
$ws.Range("O6").Value = 0.9699521281096917
$ws.Range("P6").Value = 0.9699521281096917
$ws.Range("Q6").Value = 34.64757311107711
$ws.Range("R6").Value = 311.828157999694
$ws.Range("S6").Value = 0.3014329321265276
$ws.Range("T6").Value = 0.3014329321265276

# Row 7
$ws.Range("I7").Value = 0.3107709374420163
$ws.Range("J7").Value = 0.3107709374420163
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.05133333333333333
$ws.Range("N7").Value = 0.154
$ws.Range("O7").Value = 0.00738763739326543
$ws.Range("P7").Value = 0.00738763739326543
$ws.Range("Q7").Value = 0.2638931337777778
$ws.Range("R7").Value = 2.375038204
$ws.Range("S7").Value = 0.002295862998186792
$ws.Range("T7").Value = 0.002295862998186791

# Row 8
$ws.Range("G8").Value = 2.615693
$ws.Range("H8").Value = 7.847079000000001
$ws.Range("I8").Value = 0.1581242736673807
$ws.Range("J8").Value = 0.1581242736673807
$ws.Range("O8").Value = 0.02266023449704293
$ws.Range("P8").Value = 0.02266023449704293
$ws.Range("Q8").Value = 0.4118556851103334
$ws.Range("R8").Value = 3.706701165993
$ws.Range("S8").Value = 0.003583133120977436
$ws.Range("T8").Value = 0.003583133120977436

# Row 9
$ws.Range("G9").Value = 2.615693
$ws.Range("H9").Value = 7.847079000000001
$ws.Range("I9").Value = 0.1581242736673807
$ws.Range("J9").Value = 0.1581242736673807
$ws.Range("M9").Value = 6.739756333333333
$ws.Range("O9").Value = 0.9699521281096917
$ws.Range("P9").Value = 0.9699521281096917
$ws.Range("Q9").Value = 17.62913346280567
$ws.Range("R9").Value = 158.662201165251
$ws.Range("S9").Value = 0.1533729757494751
$ws.Range("T9").Value = 0.1533729757494751

# Row 10
$ws.Range("G10").Value = 2.615693
$ws.Range("H10").Value = 7.847079000000001
$ws.Range("I10").Value = 0.1581242736673807
$ws.Range("J10").Value = 0.1581242736673807
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.05133333333333333
$ws.Range("N10").Value = 0.154
$ws.Range("O10").Value = 0.00738763739326543
$ws.Range("P10").Value = 0.00738763739326543
$ws.Range("Q10").Value = 0.1342722406666667
$ws.Range("R10").Value = 1.208450166
$ws.Range("S10").Value = 0.001168164796928078
$ws.Range("T10").Value = 0.001168164796928078

# Row 11
$ws.Range("G11").Value = 4.248598333333334
$ws.Range("H11").Value = 12.745795
$ws.Range("I11").Value = 0.2568369168563656
$ws.Range("J11").Value = 0.2568369168563656
$ws.Range("O11").Value = 0.02266023449704293
$ws.Range("P11").Value = 0.02266023449704293
$ws.Range("Q11").Value = 0.6689658829738889
$ws.Range("R11").Value = 6.020692946765
$ws.Range("S11").Value = 0.005819984763462762
$ws.Range("T11").Value = 0.005819984763462761

# Row 12
$ws.Range("G12").Value = 4.248598333333334
$ws.Range("H12").Value = 12.745795
$ws.Range("I12").Value = 0.2568369168563656
$ws.Range("J12").Value = 0.2568369168563656
$ws.Range("M12").Value = 6.739756333333333
$ws.Range("O12").Value = 0.9699521281096917
$ws.Range("P12").Value = 0.9699521281096917
$ws.Range("Q12").Value = 28.63451752487278
$ws.Range("R12").Value = 257.710657723855
$ws.Range("S12").Value = 0.2491195140819637
$ws.Range("T12").Value = 0.2491195140819637

# Row 13
$ws.Range("G13").Value = 4.248598333333334
$ws.Range("H13").Value = 12.745795
$ws.Range("I13").Value = 0.2568369168563656
$ws.Range("J13").Value = 0.2568369168563656
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.05133333333333333
$ws.Range("N13").Value = 0.154
$ws.Range("O13").Value = 0.00738763739326543
$ws.Range("P13").Value = 0.00738763739326543
$ws.Range("Q13").Value = 0.2180947144444445
$ws.Range("R13").Value = 1.96285243
$ws.Range("S13").Value = 0.001897418010939091
$ws.Range("T13").Value = 0.00189741801093909

# Row 14
$ws.Range("G14").Value = 0.9821143333333332
$ws.Range("H14").Value = 2.946343
$ws.Range("I14").Value = 0.0593709260286498
$ws.Range("J14").Value = 0.0593709260286498
$ws.Range("O14").Value = 0.02266023449704293
$ws.Range("P14").Value = 0.02266023449704293
$ws.Range("Q14").Value = 0.1546394670978889
$ws.Range("R14").Value = 1.391755203881
$ws.Range("S14").Value = 0.001345359106115795
$ws.Range("T14").Value = 0.001345359106115794

# Row 15
$ws.Range("G15").Value = 0.9821143333333332
$ws.Range("H15").Value = 2.946343
$ws.Range("I15").Value = 0.0593709260286498
$ws.Range("J15").Value = 0.0593709260286498
$ws.Range("M15").Value = 6.739756333333333
$ws.Range("O15").Value = 0.9699521281096917
$ws.Range("P15").Value = 0.9699521281096917
$ws.Range("Q15").Value = 6.619211298140777
$ws.Range("R15").Value = 59.572901683267
$ws.Range("S15").Value = 0.05758695604933196
$ws.Range("T15").Value = 0.05758695604933196

# Row 16
$ws.Range("G16").Value = 0.9821143333333332
$ws.Range("H16").Value = 2.946343
$ws.Range("I16").Value = 0.0593709260286498
$ws.Range("J16").Value = 0.0593709260286498
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.05133333333333333
$ws.Range("N16").Value = 0.154
$ws.Range("O16").Value = 0.00738763739326543
$ws.Range("P16").Value = 0.00738763739326543
$ws.Range("Q16").Value = 0.05041520244444444
$ws.Range("R16").Value = 0.4537368219999999
$ws.Range("S16").Value = 0.0004386108732020491
$ws.Range("T16").Value = 0.0004386108732020491
